$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Farmacias")

# Capture the existing hyperlink target URLs before we shift rows, so we can
# re-create them at their new location afterwards.
$wentUrl = "https://www.google.com/search?client=safari&rls=en&q=farmacia+rp+went&ie=UTF-8&oe=UTF-8"
$viamonteUrl = "https://www.google.com/search?q=farmacia+rp+viamonte&client=safari&hs=Xtj9&sca_esv=207ca25c77663f86&rls=en&sxsrf=AE3TifM4DZ0u4Ubf3LGwkNHa6GDtQ1Y8sw%3A1767101236673&ei=NNNTabjxKPjU5OUPtNewyAw&ved=0ahUKEwi42_mLteWRAxV4KrkGHbQrDMkQ4dUDCBE&uact=5&oq=farmacia+rp+viamonte&gs_lp=Egxnd3Mtd2l6LXNlcnAiFGZhcm1hY2lhIHJwIHZpYW1vbnRlMgsQLhiABBjHARivATIFEAAYgAQyBhAAGBYYHjICECYyCBAAGIAEGKIEMhoQLhiABBjHARivARiXBRjcBBjeBBjgBNgBAUjTDVC9BVicDHABeAGQAQCYAbgBoAH6B6oBAzIuNrgBA8gBAPgBAZgCCaACtwjCAgcQIxiwAxgnwgIKEAAYRxjWBBiwA8ICExAuGIAEGIoFGEMYxwEY0QMYsAOYAwCIBgGQBgq6BgYIARABGBSSBwMzLjagB9lAsgcDMi42uAexCMIHBTItOC4xyActgAgB&sclient=gws-wiz-serp"

# Insert a brand-new row above row 5 ("Farmacia Rp./ Recoleta"), pushing all
# the following rows (old 5-15) down one position (new 6-16).
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the "Farmacia Azul" record.
$ws.Range("A5").Value() = "Farmacia Azul"
$ws.Range("B5").Value() = "Av. Entre Rios 299"
$ws.Range("C5").Value() = "Balvanera"
$ws.Range("D5").Value() = "CABA"
$ws.Range("E5").Value() = "CABA"
$ws.Range("H5").Value() = 1139854640
$ws.Range("I5").Value() = "(11) 3985-3640"

# Grow the table / autofilter so the newly inserted row (and the row that
# fell off the bottom of the range) are included again: A1:J15 -> A1:J16.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J16"))

# The row-insert operation does not automatically relocate the worksheet's
# hyperlinks, so recreate them pointing at their new cells: H6->H7 (Went)
# and H7->H8 (Viamonte).
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H7"), $wentUrl, "", "", $wentUrl)
$ws.Hyperlinks.Add($ws.Range("H8"), $viamonteUrl, "", "", $viamonteUrl)

# Match the saved selection/active cell shown in the target workbook.
$ws.Activate()
$ws.Range("H5").Select()
